# "changes in the login page"
# Turns the blank Sheet1 into a small table of (username, password) test
# credentials used to validate the "invalid login" scenario of the page
# under test, with the password column hyperlinked (as Excel does
# automatically for "@"-containing text such as these password values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to describe the test data it now holds.
$ws.Name = "InvalidLoginCredentials"

# Username / password pairs.
$ws.Range("A1").Value = "anas123"
$ws.Range("B1").Value = "Zehra@2013"
$ws.Range("A2").Value = "anas234"
$ws.Range("B2").Value = "Zehra@2014"
$ws.Range("A3").Value = "anas456"
$ws.Range("B3").Value = "zehra@123"

# The password values contain "@", so Excel auto-hyperlinks them
# (mailto:) as soon as they're entered.
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:Zehra@2013") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Zehra@2014") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:zehra@123") | Out-Null

# Widen column B so the hyperlinked password text isn't clipped.
$ws.Columns.Item(2).ColumnWidth = 15.85

# Leave the selection on the row below the entered data, like after
# finishing manual data entry.
$ws.Range("A4").Select() | Out-Null
